$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.610.97'
$ws.Range("E2").Value = '  -1.34%  '
$ws.Range("D3").Value = '2.372.29'
$ws.Range("E3").Value = '  -4.31%  '
$ws.Range("E4").Value = '  -0.17%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '477.40'
$ws.Range("E5").Value = '  -2.28%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '147.09'
$ws.Range("E6").Value = '  +0.69%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.00'
$ws.Range("E7").Value = '  +0.22%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.503'
$ws.Range("E8").Value = '  -2.09%  '
$ws.Range("D9").Value = '2.374.46'
$ws.Range("E9").Value = '  -4.91%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0971'
$ws.Range("E10").Value = '  -0.08%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.41'
$ws.Range("E11").Value = '  -6.19%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.322'
$ws.Range("E12").Value = '  -2.70%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.124'
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '2.786.24'
$ws.Range("D15").Value = '55.731.47'
$ws.Range("E15").Value = '  -1.11%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '20.30'
$ws.Range("E16").Value = '  -4.42%  '
$ws.Range("E17").Value = '  -3.42%  '
$ws.Range("D18").Value = '2.376.58'
$ws.Range("E18").Value = '  -4.95%  '
$ws.Range("E19").Value = '  +0.65%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '314.94'
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("E21").Value = '  -4.74%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.999'
$ws.Range("E22").Value = '  +0.23%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.65'
$ws.Range("E23").Value = '  -2.85%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '56.78'
$ws.Range("E24").Value = '  -3.16%  '
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("E26").Value = '  -3.94%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.155'
$ws.Range("E27").Value = '  -6.05%  '
$ws.Range("D28").Value = '2.482.67'
$ws.Range("E28").Value = '  -4.65%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.23'
$ws.Range("E29").Value = '  -6.56%  '
$ws.Range("D30").Value = '0.0₃0769'
$ws.Range("E30").Value = '  -2.63%  '
$ws.Range("E31").Value = '  -0.06%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '147.92'
$ws.Range("E32").Value = '  -0.81%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '17.98'
$ws.Range("E33").Value = '  -1.45%  '
$ws.Range("E34").Value = '  -1.49%  '
$ws.Range("E35").Value = '  -2.77%  '
$ws.Range("E36").Value = '  -4.14%  '
$ws.Range("E37").Value = '  -4.30%  '
$ws.Range("E38").Value = '  -3.42%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '33.45'
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("E42").Value = '  -3.99%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.0536'
$ws.Range("E43").Value = '  -3.87%  '
$ws.Range("E44").Value = '  +3.85%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.581'
$ws.Range("E45").Value = '  -5.56%  '
$ws.Range("E46").Value = '  +0.13%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '256.16'
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("E48").Value = '  -2.24%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '4.52'
$ws.Range("E49").Value = '  -6.67%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '16.90'
$ws.Range("E50").Value = '  -4.00%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.637'
$ws.Range("E51").Value = '  +6.87%  '
